$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E2").Value = "2016-03-23 04:49:51"
$wsZh.Range("H2").Value = "2016-03-23 04:50:45"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E2").Value = "2016-03-23 04:50:00"
$wsDe.Range("H2").Value = "2016-03-23 04:50:59"
